$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Change Management Overview" ---
$ws1 = $wb.Worksheets.Item("Change Management Overview")

$ws1.Range("A2").Value = "Product Development Implementation Project"
$ws1.Range("B6").Value = "Enterprise Product Development Implementation"
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new Product Development systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in Product Development technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for Product Development transformation"

# Insert blank rows (13 and 21) that did not previously exist in the sheet,
# without disturbing the surrounding rows/data.
$ws1.Rows.Item(4).Hidden = $True
$ws1.Rows.Item(4).Hidden = $False
$ws1.Rows.Item(13).Hidden = $True
$ws1.Rows.Item(13).Hidden = $False
$ws1.Rows.Item(21).Hidden = $True
$ws1.Rows.Item(21).Hidden = $False

# --- Sheet 2: "Change Impact Assessment" ---
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")

$ws2.Range("G4").Value = "Product Development automation"
$ws2.Range("G5").Value = "Product-powered insights"
$ws2.Range("G7").Value = "New Product interface"
$ws2.Range("G11").Value = "Product-enhanced CRM"
$ws2.Range("G12").Value = "Product-assisted support"
$ws2.Range("G13").Value = "Product-powered testing"

# Insert a blank row 2 that did not previously exist in the sheet.
$ws2.Rows.Item(2).Hidden = $True
$ws2.Rows.Item(2).Hidden = $False

# --- Sheet 3: "Change Activities" ---
$ws3 = $wb.Worksheets.Item("Change Activities")

# Insert a blank row 2 that did not previously exist in the sheet.
$ws3.Rows.Item(2).Hidden = $True
$ws3.Rows.Item(2).Hidden = $False
